$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 8

$ws.Cells.Item($row, 1).Value = 42612.89402777778
$ws.Cells.Item($row, 2).Value = 22
$ws.Cells.Item($row, 3).Value = 56
$ws.Cells.Item($row, 4).Value = 43
$ws.Cells.Item($row, 5).Value = 78
$ws.Cells.Item($row, 6).Value = 21
$ws.Cells.Item($row, 7).Value = 11412
$ws.Cells.Item($row, 8).Value = 22075
$ws.Cells.Item($row, 9).Value = 2405
$ws.Cells.Item($row, 10).Value = 332
$ws.Cells.Item($row, 11).Value = 254
$ws.Cells.Item($row, 12).Value = 30
$ws.Cells.Item($row, 13).Value = 8
$ws.Cells.Item($row, 14).Value = "Bag"
